# The document has a "first page" header/footer pair plus a default
# (non-first) footer, each carrying a small logo picture as an inline
# shape:
#   - Footers(1)  -> default footer  -> Pearson logo, docPr/cNvPr id="2"
#   - Footers(2)  -> first-page footer -> Pearson logo, docPr/cNvPr id="3"
#   - Headers(2)  -> first-page header -> BTec logo,   docPr/cNvPr id="1"
#
# The authorised edit simply renames the inline pictures:
#   Pearson logo pictures: image1.png -> image2.png
#   BTec logo picture:     image2.jpg -> image1.jpg
#
# Word stores the picture's display name twice per drawing (the
# <wp:docPr name="…"/> on the drawing wrapper and the mirrored
# <pic:cNvPr name="…"/> inside the picture's non-visual properties).
# Renaming is exposed on the InlineShape object; selecting the shape
# first makes the rename stick reliably for every story (header and
# footer alike).

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlineLogo($range, $newName) {
    $shape = $range.InlineShapes(1)
    $shape.Select()
    $word.Selection.InlineShapes(1).Name = $newName
}

# Pearson logo in the default footer (footer2.xml, id="2")
Rename-InlineLogo $sec.Footers(1).Range "image2.png"

# Pearson logo in the first-page footer (footer1.xml, id="3")
Rename-InlineLogo $sec.Footers(2).Range "image2.png"

# BTec logo in the first-page header (header1.xml, id="1")
Rename-InlineLogo $sec.Headers(2).Range "image1.jpg"

Write-Output "Renamed Pearson logo pictures to image2.png and BTec logo picture to image1.jpg"
